$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 507
$ws1.Range("F7").Value = 707

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 507
$ws4.Range("F7").Value = 707
